$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 400.33334
$ws.Range("I11").Value = 400.33334
$ws.Range("K11").Value = 400.33334
$ws.Range("M11").Value = -260.33334
# Row 17
$ws.Range("H17").Value = 2735.3333
$ws.Range("J17").Value = 2420.9092
$ws.Range("L17").Value = 7262.7276
$ws.Range("N17").Value = -7598.7276
# Row 19
$ws.Range("H19").Value = 430.75
$ws.Range("I19").Value = 300
$ws.Range("J19").Value = 474.33334
$ws.Range("K19").Value = 300
$ws.Range("L19").Value = 474.33334
$ws.Range("M19").Value = -125
$ws.Range("N19").Value = -824.33334
# Row 40
$ws.Range("H40").Value = 5849.3335
$ws.Range("I40").Value = 5424.5
$ws.Range("J40").Value = 6699
$ws.Range("K40").Value = 5424.5
$ws.Range("L40").Value = 6699
$ws.Range("M40").Value = -5249.5
$ws.Range("N40").Value = -7049
# Row 70
$ws.Range("H70").Value = 850
$ws.Range("I70").Value = 525
$ws.Range("J70").Value = 1500
$ws.Range("K70").Value = 1575
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = -1305
$ws.Range("N70").Value = -5040
# Row 73
$ws.Range("H73").Value = 850
$ws.Range("I73").Value = 525
$ws.Range("J73").Value = 1500
$ws.Range("K73").Value = 1575
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = -639
$ws.Range("N73").Value = -6372
# Row 87
$ws.Range("H87").Value = 103329
$ws.Range("J87").Value = 103329
$ws.Range("L87").Value = 103329
$ws.Range("N87").Value = -105825
# Row 90
$ws.Range("H90").Value = 103329
$ws.Range("J90").Value = 103329
$ws.Range("L90").Value = 309987
$ws.Range("N90").Value = -322467
# Row 97
$ws.Range("H97").Value = 9275.286
$ws.Range("J97").Value = 9275.286
$ws.Range("L97").Value = 27825.858
$ws.Range("N97").Value = -28817.858
# Row 138
$ws.Range("H138").Value = 3485.375
$ws.Range("J138").Value = 4121.5
$ws.Range("L138").Value = 12364.5
$ws.Range("N138").Value = -22644.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3665.276
$ws.Range("I32").Value = 3581.9285
$ws.Range("K32").Value = 3581.9285
$ws.Range("M32").Value = -3294.9285
# Row 61
$ws.Range("H61").Value = 2229.558
$ws.Range("I61").Value = 1642.9706
$ws.Range("J61").Value = 4445.5557
$ws.Range("K61").Value = 1642.9706
$ws.Range("L61").Value = 4445.5557
$ws.Range("M61").Value = -1430.9706
$ws.Range("N61").Value = -4869.5557
# Row 63
$ws.Range("H63").Value = 2074.8572
$ws.Range("I63").Value = 2213.4546
$ws.Range("J63").Value = 1566.6666
$ws.Range("K63").Value = 2213.4546
$ws.Range("L63").Value = 1566.6666
$ws.Range("M63").Value = -1527.4546
$ws.Range("N63").Value = -2938.6666
# Row 66
$ws.Range("H66").Value = 2074.8572
$ws.Range("I66").Value = 2213.4546
$ws.Range("J66").Value = 1566.6666
$ws.Range("K66").Value = 11067.273
$ws.Range("L66").Value = 7833.333000000001
$ws.Range("M66").Value = -7635.273000000001
$ws.Range("N66").Value = -14697.333
# Row 132
$ws.Range("H132").Value = 3049.76
$ws.Range("I132").Value = 2724.825
$ws.Range("K132").Value = 8174.474999999999
$ws.Range("M132").Value = -5644.474999999999
# Row 136
$ws.Range("H136").Value = 2229.558
$ws.Range("I136").Value = 1642.9706
$ws.Range("J136").Value = 4445.5557
$ws.Range("K136").Value = 4928.9118
$ws.Range("L136").Value = 13336.6671
$ws.Range("M136").Value = -2378.9118
$ws.Range("N136").Value = -18436.6671

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 353.25
$ws.Range("I22").Value = 419.8
$ws.Range("J22").Value = 242.33333
$ws.Range("K22").Value = 419.8
$ws.Range("L22").Value = 242.33333
$ws.Range("M22").Value = -246.8
$ws.Range("N22").Value = -588.3333299999999
# Row 99
$ws.Range("H99").Value = 4108.0713
$ws.Range("I99").Value = 2703.1052
$ws.Range("K99").Value = 2703.1052
$ws.Range("M99").Value = -1205.1052
# Row 134
$ws.Range("H134").Value = 3172.5789
$ws.Range("I134").Value = 2865.394
$ws.Range("K134").Value = 8596.181999999999
$ws.Range("M134").Value = -6061.181999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2737
$ws.Range("I58").Value = 1403
$ws.Range("K58").Value = 1403
$ws.Range("M58").Value = -1200
# Row 86
$ws.Range("H86").Value = 23361
$ws.Range("J86").Value = 8224.25
$ws.Range("L86").Value = 8224.25
$ws.Range("N86").Value = -10470.25
# Row 89
$ws.Range("H89").Value = 23361
$ws.Range("J89").Value = 8224.25
$ws.Range("L89").Value = 41121.25
$ws.Range("N89").Value = -52353.25
# Row 99
$ws.Range("H99").Value = 3996.3
$ws.Range("I99").Value = 3990.389
$ws.Range("K99").Value = 3990.389
$ws.Range("M99").Value = -2492.389
# Row 126
$ws.Range("H126").Value = 3996.3
$ws.Range("I126").Value = 3990.389
$ws.Range("K126").Value = 11971.167
$ws.Range("M126").Value = -9501.167000000001
# Row 132
$ws.Range("H132").Value = 50002580
$ws.Range("I132").Value = 71429976
$ws.Range("J132").Value = 5333
$ws.Range("K132").Value = 214289928
$ws.Range("L132").Value = 15999
$ws.Range("M132").Value = -214287398
$ws.Range("N132").Value = -21059
# Row 134
$ws.Range("H134").Value = 3637.7856
$ws.Range("I134").Value = 3542.9
$ws.Range("K134").Value = 10628.7
$ws.Range("M134").Value = -8093.700000000001
# Row 136
$ws.Range("H136").Value = 2737
$ws.Range("I136").Value = 1403
$ws.Range("K136").Value = 4209
$ws.Range("M136").Value = -1659

$ws = $wb.Worksheets.Item("CUL")
# Row 59
$ws.Range("H59").Value = 500
$ws.Range("I59").Value = 500
$ws.Range("K59").Value = 1500
$ws.Range("M59").Value = -960
# Row 97
$ws.Range("H97").Value = 743321.3
$ws.Range("I97").Value = 1666999.6
$ws.Range("J97").Value = 50562.5
$ws.Range("K97").Value = 5000998.800000001
$ws.Range("L97").Value = 151687.5
$ws.Range("M97").Value = -5000502.800000001
$ws.Range("N97").Value = -152679.5
# Row 107
$ws.Range("H107").Value = 378.4
$ws.Range("J107").Value = 378.4
$ws.Range("L107").Value = 1135.2
$ws.Range("N107").Value = -4975.2

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1107.2727
$ws.Range("I102").Value = 1080.4
$ws.Range("K102").Value = 1080.4
$ws.Range("M102").Value = 541.5999999999999
# Row 132
$ws.Range("H132").Value = 2518.7234
$ws.Range("I132").Value = 2017.8948
$ws.Range("K132").Value = 6053.6844
$ws.Range("M132").Value = -3523.6844

$ws = $wb.Worksheets.Item("LTW")
# Row 36
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
# Row 46
$ws.Range("H46").Value = 2516.16
$ws.Range("I46").Value = 2483
$ws.Range("J46").Value = 2586.625
$ws.Range("K46").Value = 2483
$ws.Range("L46").Value = 2586.625
$ws.Range("M46").Value = -2295
$ws.Range("N46").Value = -2962.625
# Row 93
$ws.Range("H93").Value = 1294.3334
$ws.Range("I93").Value = 1439.5
$ws.Range("J93").Value = 1004
$ws.Range("K93").Value = 1439.5
$ws.Range("L93").Value = 1004
$ws.Range("M93").Value = -191.5
$ws.Range("N93").Value = -3500
# Row 132
$ws.Range("H132").Value = 4521.5
$ws.Range("I132").Value = 2263.3333
$ws.Range("K132").Value = 6789.999899999999
$ws.Range("M132").Value = -4259.999899999999
# Row 136
$ws.Range("H136").Value = 4496.826
$ws.Range("I136").Value = 4012.6428
$ws.Range("K136").Value = 12037.9284
$ws.Range("M136").Value = -9487.928400000001

$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
# Row 100
$ws.Range("H100").Value = 111112150
$ws.Range("I100").Value = 899.2
$ws.Range("K100").Value = 1798.4
$ws.Range("M100").Value = -1257.4
# Row 132
$ws.Range("H132").Value = 19612732
$ws.Range("I132").Value = 25644840
$ws.Range("J132").Value = 8374.5
$ws.Range("K132").Value = 76934520
$ws.Range("L132").Value = 25123.5
$ws.Range("M132").Value = -76931990
$ws.Range("N132").Value = -30183.5
# Row 134
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
# Row 136
$ws.Range("H136").Value = 22216.87
$ws.Range("J136").Value = 27901
$ws.Range("L136").Value = 83703
$ws.Range("N136").Value = -88803

